$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("icf_a.185")
$ws.Activate()

# Update row 12 values from BW..CG (columns 75..85) to 1540
$ws.Range("BW12:CG12").Value = 1540

# Move selection to CG13 to match the saved view state
$ws.Range("CG13").Select()
